# Transition rule summary tables: add "Within 5 miles" and "Within 10 miles"
# of HFC production facility columns (F, G) to both the Means sheet and the
# Standard Deviations sheet, and update a handful of existing values that
# were recomputed alongside the new radii (Total Cancer Risk / Total
# Respiratory rows on both sheets, plus their B/C columns).

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: "Means"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New header cells for the two new distance bands
$ws1.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws1.Range("G1").Value = "Within 10 miles of HFC production facility"

# Row 2: % White
$ws1.Range("F2").Value = 79
$ws1.Range("G2").Value = 71

# Row 3: % Black or African American
$ws1.Range("F3").Value = 18
$ws1.Range("G3").Value = 27

# Row 4: % Other
$ws1.Range("F4").Value = 2.7
$ws1.Range("G4").Value = 2.2

# Row 5: % Hispanic
$ws1.Range("F5").Value = 3.5
$ws1.Range("G5").Value = 1.9

# Row 6: Median Income [1,000 2019$]
$ws1.Range("F6").Value = 39
$ws1.Range("G6").Value = 38

# Row 7: % Below Poverty Line
$ws1.Range("F7").Value = 14
$ws1.Range("G7").Value = 15

# Row 8: % Below Half the Poverty Line
$ws1.Range("F8").Value = 10
$ws1.Range("G8").Value = 8.9

# Row 9: Total Cancer Risk (per million) - B/C recomputed, plus new F/G
$ws1.Range("B9").Value = 26
$ws1.Range("C9").Value = 39
$ws1.Range("F9").Value = 40
$ws1.Range("G9").Value = 40

# Row 10: Total Respiratory (hazard quotient) - B/C/D/E recomputed, plus new F/G
$ws1.Range("B10").Value = 0.32
$ws1.Range("C10").Value = 0.43
$ws1.Range("D10").Value = 0.5
$ws1.Range("E10").Value = 0.5
$ws1.Range("F10").Value = 0.5
$ws1.Range("G10").Value = 0.5

# ----------------------------------------------------------------------
# Sheet 2: "Standard Deviations"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# New header cells for the two new distance bands
$ws2.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws2.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# Row 2: % White
$ws2.Range("F2").Value = 15
$ws2.Range("G2").Value = 17

# Row 3: % Black or African American
$ws2.Range("F3").Value = 15
$ws2.Range("G3").Value = 18

# Row 4: % Other
$ws2.Range("F4").Value = 2.4
$ws2.Range("G4").Value = 2

# Row 5: % Hispanic
$ws2.Range("F5").Value = 4
$ws2.Range("G5").Value = 3

# Row 6: Median Income [1,000 2019$]
$ws2.Range("F6").Value = 13
$ws2.Range("G6").Value = 14

# Row 7: % Below Poverty Line
$ws2.Range("F7").Value = 9.1
$ws2.Range("G7").Value = 11

# Row 8: % Below Half the Poverty Line
$ws2.Range("F8").Value = 7.4
$ws2.Range("G8").Value = 6.6

# Row 9: Total Cancer Risk (per million) - B/C recomputed, D/E now 0, plus new F/G
$ws2.Range("B9").Value = 8.6
$ws2.Range("C9").Value = 24
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 0
$ws2.Range("F9").Value = 0
$ws2.Range("G9").Value = 0

# Row 10: Total Respiratory (hazard quotient) - C recomputed, D/E now 0, plus new F/G
$ws2.Range("B10").Value = 0.14
$ws2.Range("C10").Value = 0.084
$ws2.Range("D10").Value = 0
$ws2.Range("E10").Value = 0
$ws2.Range("F10").Value = 0
$ws2.Range("G10").Value = 0
